$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at row 40 (shifts old total/footer rows down) ---
$ws.Rows.Item(40).Insert()

# Clone formatting from the row above (row 39 - last existing data row)
$ws.Range("A39:Q39").Copy()
$ws.Range("A40:Q40").PasteSpecial(-4122)
$ws.Rows.Item(40).RowHeight = 24.75

# Re-create the merges for the new data row (mirroring row 39's merge layout)
$ws.Range("A40:B40").Merge()
$ws.Range("C40:G40").Merge()
$ws.Range("H40:K40").Merge()
$ws.Range("L40:M40").Merge()
$ws.Range("N40:O40").Merge()

# --- Populate the new data row with the new item ---
$ws.Range("A40").Value = 34
$ws.Range("C40").Value = "نيفيا سوفت كريم 50 مل"
$ws.Range("H40").Value = "1:0"

# L40 / P40 are numeric-formatted columns but need to hold literal text values
# (matching how the rest of the sheet stores these "numbers" as shared strings).
$ws.Range("L40").NumberFormat = "@"
$ws.Range("L40").Value = "0"
$ws.Range("L40").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N40").Value = "40.00"

$ws.Range("P40").NumberFormat = "@"
$ws.Range("P40").Value = "40.0000"
$ws.Range("P40").NumberFormat = "0.00"

$ws.Range("Q40").Value = "1:0"

# --- Update the totals row (now row 41): add the new item's price ---
$ws.Range("P41").Value = 1255.915

# --- Update the footer row (now row 42): refresh the generated timestamp ---
$ws.Range("A42").Value = "Thursday, 31 July, 2025 3:15 PM"
